$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: Price values in column D are stored as plain text in this sheet
# (e.g. "493.99"), not numbers. Assigning a numeric-looking string via
# .Value would make Excel auto-convert it to a real number, so those
# assignments are prefixed with a leading apostrophe ('') — the standard
# Excel "treat as text" quote-prefix — to keep them as text, matching the
# original inline-string cell type.

$ws.Range("D2").Value = '54.246.99'
$ws.Range("E2").Value = '  -3.20%  '

$ws.Range("D3").Value = '2.289.29'
$ws.Range("E3").Value = '  -3.31%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = '''493.99'
$ws.Range("E5").Value = '  -1.92%  '

$ws.Range("D6").Value = '''128.05'
$ws.Range("E6").Value = '  -2.32%  '

$ws.Range("E7").Value = '  -0.11%  '

$ws.Range("D8").Value = '''0.527'
$ws.Range("E8").Value = '  -3.83%  '

$ws.Range("D9").Value = '2.298.77'
$ws.Range("E9").Value = '  -2.96%  '

$ws.Range("D10").Value = '''0.0938'
$ws.Range("E10").Value = '  -3.61%  '

$ws.Range("E11").Value = '  -1.29%  '

$ws.Range("E12").Value = '  +1.66%  '

$ws.Range("E13").Value = '  -3.04%  '

$ws.Range("D14").Value = '2.697.01'

$ws.Range("D15").Value = '''21.39'
$ws.Range("E15").Value = '  -0.30%  '

$ws.Range("D16").Value = '54.226.59'
$ws.Range("E16").Value = '  -3.13%  '

$ws.Range("E17").Value = '  -1.96%  '

$ws.Range("D18").Value = '2.285.99'
$ws.Range("E18").Value = '  -4.88%  '

$ws.Range("E19").Value = '  +0.16%  '

$ws.Range("D20").Value = '''9.68'
$ws.Range("E20").Value = '  -3.38%  '

$ws.Range("D21").Value = '''303.49'
$ws.Range("E21").Value = '  -1.14%  '

$ws.Range("D22").Value = '''6.18'
$ws.Range("E22").Value = '  -1.13%  '

$ws.Range("E23").Value = '  -0.14%  '

$ws.Range("D24").Value = '''64.03'
$ws.Range("E24").Value = '  -1.61%  '

$ws.Range("D25").Value = '''0.997'
$ws.Range("E25").Value = '  -0.51%  '

$ws.Range("D26").Value = '''0.368'
$ws.Range("E26").Value = '  -0.68%  '

$ws.Range("D28").Value = '''7.11'
$ws.Range("E28").Value = '  -1.69%  '

$ws.Range("D29").Value = '''169.71'
$ws.Range("E29").Value = '  -1.79%  '

$ws.Range("D30").Value = '0.0₃0700'
$ws.Range("E30").Value = '  -2.28%  '

$ws.Range("E31").Value = '  -1.41%  '

$ws.Range("E32").Value = '  -0.11%  '

$ws.Range("E33").Value = '  +0.95%  '

$ws.Range("D34").Value = '''0.998'
$ws.Range("E34").Value = '  +0.00%  '

$ws.Range("D35").Value = '''1.07'
$ws.Range("E35").Value = '  -2.19%  '

$ws.Range("D36").Value = '''17.61'
$ws.Range("E36").Value = '  +0.14%  '

$ws.Range("E37").Value = '  -1.13%  '

$ws.Range("D38").Value = '''0.848'
$ws.Range("E38").Value = '  +6.45%  '

$ws.Range("E39").Value = '  -4.07%  '

$ws.Range("D40").Value = '''35.81'
$ws.Range("E40").Value = '  -0.59%  '

$ws.Range("D41").Value = '''1.38'
$ws.Range("E41").Value = '  -2.22%  '

$ws.Range("D42").Value = '''0.369'
$ws.Range("E42").Value = '  -0.07%  '

$ws.Range("D43").Value = '''3.33'
$ws.Range("E43").Value = '  -0.59%  '

$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").Value = '''4.80'
$ws.Range("E44").Value = '  +0.02%  '

$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = '''123.47'
$ws.Range("E45").Value = '  -5.85%  '

$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").Value = '''0.547'
$ws.Range("E46").Value = '  -2.80%  '

$ws.Range("B47").Value = 'Stellar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D47").Value = '''0.0883'
$ws.Range("E47").Value = '  -2.76%  '

$ws.Range("D48").Value = '''238.68'
$ws.Range("E48").Value = '  -2.29%  '

$ws.Range("E49").Value = '  -0.91%  '

$ws.Range("D50").Value = '''0.0204'
$ws.Range("E50").Value = '  -1.55%  '

$ws.Range("E51").Value = '  -2.18%  '
